$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend header row (row 1) with two more sequential values: P1=14, Q1=15
$ws.Cells.Item(1, 16).Value = 14
$ws.Cells.Item(1, 17).Value = 15

# Copy the style (bold, centered, thin border) from O1 onto the new header cells.
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- New repeating data pattern for columns B:Q (cols 2..17), identical for every data row 2-25.
$rowValues = @(2, 2, 2, 1, 1, 1, 2, 2, 2, 1, 2, 2, 2, 1, 2, 2)

for ($r = 2; $r -le 25; $r++) {
    for ($ci = 0; $ci -lt $rowValues.Length; $ci++) {
        $ws.Cells.Item($r, $ci + 2).Value = $rowValues[$ci]
    }
}
